$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.70"
$ws.Range("E2").Value = "'3.16%"
$ws.Range("D3").Value = "'41.44"
$ws.Range("E3").Value = "'5.88%"
$ws.Range("D4").Value = "'5.666"
$ws.Range("E4").Value = "'-4.00%"
$ws.Range("D5").Value = "'0.08184"
$ws.Range("E5").Value = "'2.22%"
$ws.Range("D6").Value = "'2.076"
$ws.Range("E6").Value = "'10.98%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.718"
$ws.Range("E7").Value = "'0.87%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.949"
$ws.Range("E8").Value = "'0.20%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9236"
$ws.Range("E9").Value = "'-1.13%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1255"
$ws.Range("E10").Value = "'1.80%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1955"
$ws.Range("E11").Value = "'0.27%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09396"
$ws.Range("E12").Value = "'2.73%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03685"
$ws.Range("E13").Value = "'7.54%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1054"
$ws.Range("E14").Value = "'9.94%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001286"
$ws.Range("E15").Value = "'-0.45%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006196"
$ws.Range("E16").Value = "'1.03%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.410"
$ws.Range("E17").Value = "'1.95%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.541"
$ws.Range("E18").Value = "'-0.88%"
$ws.Range("E19").Value = "'-1.44%"
$ws.Range("D20").Value = "'8.321"
$ws.Range("E20").Value = "'-5.02%"
$ws.Range("D21").Value = "'0.1388"
$ws.Range("E21").Value = "'-1.67%"
$ws.Range("E22").Value = "'9.92%"
$ws.Range("D23").Value = "'0.04445"
$ws.Range("E23").Value = "'-0.67%"
$ws.Range("D24").Value = "'0.001273"
$ws.Range("E24").Value = "'0.59%"
$ws.Range("D25").Value = "'0.004321"
$ws.Range("E25").Value = "'-0.81%"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'4.98%"
$ws.Range("D39").Value = "'0.02813"
$ws.Range("E39").Value = "'16.66%"
$ws.Range("D40").Value = "'0.05452"
$ws.Range("E40").Value = "'5.28%"
$ws.Range("D41").Value = "'0.007675"
$ws.Range("E41").Value = "'3.18%"
$ws.Range("D42").Value = "'0.009437"
$ws.Range("E42").Value = "'3.52%"
$ws.Range("D43").Value = "'0.1418"
$ws.Range("E43").Value = "'1.00%"
$ws.Range("D44").Value = "'0.002132"
$ws.Range("E44").Value = "'1.58%"
$ws.Range("D45").Value = "'0.01185"
$ws.Range("E45").Value = "'3.42%"
$ws.Range("D46").Value = "'0.00006869"
$ws.Range("E46").Value = "'1.69%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.14%"
$ws.Range("D48").Value = "'0.002282"
$ws.Range("E48").Value = "'60.38%"
$ws.Range("E49").Value = "'7.48%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.14%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.14%"
